# Daily attendance processing - 2025-11-20 23:43:29
#
# Column G ("Recorded By") holds a comma-separated list of the users who
# recorded/edited each attendance entry. Normalize the ordering so "System"
# is listed first whenever it appears alongside other recorders.
#
# Rule: for any multi-value "Recorded By" cell whose first entry is not
# already "System", reverse the order of the comma-separated entries.
# Single-value cells, and cells that already start with "System", are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = $ws.UsedRange.Rows.Count
$col = 7   # column G ("Recorded By")

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $count = $parts.Count

        if ($count -gt 1 -and $parts[0] -ne "System") {
            $reversed = @()
            for ($i = $count - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = ($reversed -join ", ")
        }
    }
}
